# Library/SmokeSanityTestCases.xlsx - refresh the smoke/sanity test case list.
# The sheet held 17 data rows (rows 2-18); the refreshed list only has
# 14 data rows (rows 2-15), with updated test-case names. Status stays
# "Norun" for every row, Plan is SmokeTest for the first 7 rows and
# SanityTest for the remaining 7, and TC No. is simply TC001..TC014.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three trailing rows that no longer exist in the refreshed list.
$ws.Rows("16:18").Delete()

# Clear any leftover notes in column F (Deprecated/Hidden/OnHold markers)
# for the rows that remain - the refreshed list carries no such notes.
$ws.Range("F2:F15").ClearContents()

# New set of test cases (TestCaseName, Status, Plan, Run No., TC No., Bug, Manual Status)
$rows = @(
    @{ Row=2;  Name="Appointment\TC001CreateAppointmentNew.py";                Plan="SmokeTest";  TC="TC001" },
    @{ Row=3;  Name="Billing\Opbilling\TC002OPDbillingLabXray.py";              Plan="SmokeTest";  TC="TC002" },
    @{ Row=4;  Name="Laboratory\TC005GenerateLabReport.py";                     Plan="SmokeTest";  TC="TC003" },
    @{ Row=5;  Name="Radiology\TC001GenerateUSGReport.py";                      Plan="SmokeTest";  TC="TC004" },
    @{ Row=6;  Name="Dispensary\TC001CreateDispensarySale.py";                  Plan="SmokeTest";  TC="TC005" },
    @{ Row=7;  Name="ADT\TC010AdmissionDischargeTransferToBePaid.py";           Plan="SmokeTest";  TC="TC006" },
    @{ Row=8;  Name="Appointment\TC002CreateAppointmentFollowup.py";            Plan="SmokeTest";  TC="TC007" },
    @{ Row=9;  Name="Reports\TC011UserCollectionReport.py";                     Plan="SanityTest"; TC="TC008" },
    @{ Row=10; Name="Reports\TC002TotalItemsBillReport.py";                     Plan="SanityTest"; TC="TC009" },
    @{ Row=11; Name="Reports\TC001BillingDashboardSummary.py";                  Plan="SanityTest"; TC="TC010" },
    @{ Row=12; Name="Pharmacy\Reports\TC010VerifyPharmacyDashboard.py";         Plan="SanityTest"; TC="TC011" },
    @{ Row=13; Name="Pharmacy\TC003PharmacyOPDbilling.py";                      Plan="SanityTest"; TC="TC012" },
    @{ Row=14; Name="Pharmacy\Reports\TC002UserCollectionReport.py";            Plan="SanityTest"; TC="TC013" },
    @{ Row=15; Name="Inventory\TC002Verify_DirectDispatch&PurchaseRequest.py";  Plan="SanityTest"; TC="TC014" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.Name
    $ws.Range("B$n").Value = "Norun"
    $ws.Range("C$n").Value = $r.Plan
    $ws.Range("E$n").Value = $r.TC
}

# The stray "J" column marker (single space) moves from row 7 to row 8
# along with the reshuffled test case rows.
$ws.Range("J7").ClearContents()
$ws.Range("J8").Value = " "

# Update the active cell/selection to reflect the trimmed sheet.
$ws.Range("E20").Select()

Write-Host "Library smoke/sanity test case list refreshed."
